$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Devices")

# --- New column J: "DC Unit Loading Details Name" header + "Current (DC Units)" value ---

# Header cell J7: copy the look of the existing header cell I7 (bold font, blue fill,
# full thin border) then strip the top/bottom border so only left/right remain,
# matching the new distinct header style used for this column.
$ws.Range("I7").Copy()
$ws.Range("J7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J7").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none
$ws.Range("J7").Borders.Item(8).LineStyle = -4142   # xlEdgeTop -> none
$ws.Range("J7").Value = "DC Unit Loading Details Name"

# Data cell J8: copy the look of an existing data cell (A8) which uses the
# standard body style (thin border, light fill, wrapped + vertically centered).
$ws.Range("A8").Copy()
$ws.Range("J8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J8").Value = "Current (DC Units)"

# Widen column J to fit the new, longer header text.
$ws.Columns.Item(10).ColumnWidth = 25.5

# Leave the final selection on the last cell that was edited.
$null = $ws.Range("J8").Select()

$excel.CutCopyMode = 0
